$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18; this shifts the existing rows 18-51
# down to 19-52, preserving all of their data.
$ws.Rows("18").Insert()

# Populate the newly inserted row 18 with the new price-report entry.
$ws.Range("A18").Value = 10
$ws.Range("B18").Value = "Vega Modelo de Temuco"
$ws.Range("C18").Value = "La Araucanía"
$ws.Range("D18").Value = 44526
$ws.Range("E18").Value = 9
$ws.Range("F18").Value = 300000001
$ws.Range("G18").Value = "Rabanito"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 50
$ws.Range("K18").Value = 6000
$ws.Range("L18").Value = 7000
$ws.Range("M18").Value = 6400
$ws.Range("N18").Value = "$/docena de paquetes"
$ws.Range("O18").Value = "Provincia de Cautín"
$ws.Range("P18").Value = 533
$ws.Range("Q18").Value = 12
$ws.Range("R18").Value = "Hortaliza"
